# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates metric values in row 3 (metrics_sim_with_priors.json) of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.6486486486486487
$ws.Range("D3").Value = 0.9459459459459459

$ws.Range("H3").Value = 0.6870748299319728
$ws.Range("I3").Value = 0.07570766032304493
$ws.Range("J3").Value = 0.5405405405405406
$ws.Range("K3").Value = 82.97297297297297

$ws.Range("Q3").Value = 8
$ws.Range("R3").Value = 17
$ws.Range("S3").Value = 38
$ws.Range("T3").Value = 79
$ws.Range("U3").Value = 161

$ws.Range("V3").Value = 837
$ws.Range("W3").Value = 828
$ws.Range("X3").Value = 807
$ws.Range("Y3").Value = 766
$ws.Range("Z3").Value = 684

$ws.Range("AF3").Value = 0.990533
$ws.Range("AG3").Value = 0.979882
$ws.Range("AH3").Value = 0.95503
$ws.Range("AI3").Value = 0.906509
$ws.Range("AJ3").Value = 0.809467
